# 11/30/2015 add new object in member account
#
# Updates TestData.xlsx / Sheet1:
#   - Rename header cells (A1: "Test Case Name" -> "Test_Name",
#     E1: "Actuel SAndBox" -> "Actuel Result")
#   - Rename the test-case id in A2 ("TS 1.1" -> "1.TC.1.1")
#   - Add the new "member account" object column: E2:E5 = 1
#   - Remove the now-unused scratch rows 10:13
#   - Resize columns A and E, add widths for the (new) F and H columns
#   - Move the active selection to E10

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Header row -------------------------------------------------------
$ws.Range("A1").Value = "Test_Name"

# --- Test case row 2 ----------------------------------------------------
$ws.Range("A2").Value = "1.TC.1.1"

$ws.Range("E1").Value = "Actuel Result"

# --- New "member account" object values (column E) ----------------------
$ws.Range("E2").Value = 1
$ws.Range("E3").Value = 1
$ws.Range("E4").Value = 1
$ws.Range("E5").Value = 1

# --- Remove the old scratch/debug rows -----------------------------------
$ws.Rows("10:13").Delete()

# --- Column width adjustments --------------------------------------------
# (ColumnWidth values chosen so the saved <col width="..."> lands as close
#  as possible to the target 22.42578125 / 15.85546875 / 20.28515625 / 18)
$ws.Columns.Item(1).ColumnWidth = 21.666666666666668
$ws.Columns.Item(5).ColumnWidth = 15.0
$ws.Columns.Item(6).ColumnWidth = 19.5
$ws.Columns.Item(8).ColumnWidth = 17.166666666666668

# --- Selection -------------------------------------------------------------
$ws.Range("E10").Select() | Out-Null
